$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(28, 8).Value2 = 568.875  # H28: 522.2222 -> 568.875
$ws.Cells.Item(28, 9).Value2 = 593  # I28: 537.5 -> 593
$ws.Cells.Item(28, 11).Value2 = 593  # K28: 537.5 -> 593
$ws.Cells.Item(28, 13).Value2 = -108  # M28: -52.5 -> -108

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(69, 8).Value2 = 8178.25  # H69: 9362 -> 8178.25
$ws.Cells.Item(69, 9).Value2 = 5015  # I69: 0 -> 5015
$ws.Cells.Item(69, 10).Value2 = 9232.666999999999  # J69: 9362 -> 9232.666999999999
$ws.Cells.Item(69, 11).Value2 = 15045  # K69: 0 -> 15045
$ws.Cells.Item(69, 12).Value2 = 27698.001  # L69: 28086 -> 27698.001
$ws.Cells.Item(69, 13).Value2 = -14171  # M69: None -> -14171
$ws.Cells.Item(69, 14).Value2 = -29446.001  # N69: -29834 -> -29446.001

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(70, 8).Value2 = 86545.89  # H70: 12904.667 -> 86545.89
$ws.Cells.Item(70, 9).Value2 = 1185.5  # I70: 880 -> 1185.5
$ws.Cells.Item(70, 10).Value2 = 154834.2  # J70: 18917 -> 154834.2
$ws.Cells.Item(70, 11).Value2 = 3556.5  # K70: 2640 -> 3556.5
$ws.Cells.Item(70, 12).Value2 = 464502.6  # L70: 56751 -> 464502.6
$ws.Cells.Item(70, 13).Value2 = -3286.5  # M70: -2370 -> -3286.5
$ws.Cells.Item(70, 14).Value2 = -465042.6  # N70: -57291 -> -465042.6

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(72, 8).Value2 = 8178.25  # H72: 9362 -> 8178.25
$ws.Cells.Item(72, 9).Value2 = 5015  # I72: 0 -> 5015
$ws.Cells.Item(72, 10).Value2 = 9232.666999999999  # J72: 9362 -> 9232.666999999999
$ws.Cells.Item(72, 11).Value2 = 45135  # K72: 0 -> 45135
$ws.Cells.Item(72, 12).Value2 = 83094.003  # L72: 84258 -> 83094.003
$ws.Cells.Item(72, 13).Value2 = -40767  # M72: None -> -40767
$ws.Cells.Item(72, 14).Value2 = -91830.003  # N72: -92994 -> -91830.003

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(73, 8).Value2 = 86545.89  # H73: 12904.667 -> 86545.89
$ws.Cells.Item(73, 9).Value2 = 1185.5  # I73: 880 -> 1185.5
$ws.Cells.Item(73, 10).Value2 = 154834.2  # J73: 18917 -> 154834.2
$ws.Cells.Item(73, 11).Value2 = 3556.5  # K73: 2640 -> 3556.5
$ws.Cells.Item(73, 12).Value2 = 464502.6  # L73: 56751 -> 464502.6
$ws.Cells.Item(73, 13).Value2 = -2620.5  # M73: -1704 -> -2620.5
$ws.Cells.Item(73, 14).Value2 = -466374.6  # N73: -58623 -> -466374.6

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(98, 8).Value2 = 1100.1818  # H98: 1150.7 -> 1100.1818
$ws.Cells.Item(98, 9).Value2 = 1100.1818  # I98: 1150.7 -> 1100.1818
$ws.Cells.Item(98, 11).Value2 = 1100.1818  # K98: 1150.7 -> 1100.1818
$ws.Cells.Item(98, 13).Value2 = 397.8181999999999  # M98: 347.3 -> 397.8181999999999

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(122, 8).Value2 = 1100.1818  # H122: 1150.7 -> 1100.1818
$ws.Cells.Item(122, 9).Value2 = 1100.1818  # I122: 1150.7 -> 1100.1818
$ws.Cells.Item(122, 11).Value2 = 3300.5454  # K122: 3452.1 -> 3300.5454
$ws.Cells.Item(122, 13).Value2 = -850.5454  # M122: -1002.1 -> -850.5454

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value2 = 2389.2144  # H132: 2803.6365 -> 2389.2144
$ws.Cells.Item(132, 9).Value2 = 2389.2144  # I132: 2803.6365 -> 2389.2144
$ws.Cells.Item(132, 11).Value2 = 7167.6432  # K132: 8410.9095 -> 7167.6432
$ws.Cells.Item(132, 13).Value2 = -4637.6432  # M132: -5880.9095 -> -4637.6432

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value2 = 1804.6552  # H137: 2034.5172 -> 1804.6552
$ws.Cells.Item(137, 9).Value2 = 952.44446  # I137: 1077.0588 -> 952.44446
$ws.Cells.Item(137, 10).Value2 = 3199.182  # J137: 3390.9167 -> 3199.182
$ws.Cells.Item(137, 11).Value2 = 2857.33338  # K137: 3231.1764 -> 2857.33338
$ws.Cells.Item(137, 12).Value2 = 9597.545999999998  # L137: 10172.7501 -> 9597.545999999998
$ws.Cells.Item(137, 13).Value2 = -307.33338  # M137: -681.1764000000003 -> -307.33338
$ws.Cells.Item(137, 14).Value2 = -14697.546  # N137: -15272.7501 -> -14697.546

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(141, 8).Value2 = 2910  # H141: 3472.5 -> 2910
$ws.Cells.Item(141, 9).Value2 = 2176  # I141: 2745 -> 2176
$ws.Cells.Item(141, 10).Value2 = 4133.3335  # J141: 4200 -> 4133.3335
$ws.Cells.Item(141, 11).Value2 = 6528  # K141: 8235 -> 6528
$ws.Cells.Item(141, 12).Value2 = 12400.0005  # L141: 12600 -> 12400.0005
$ws.Cells.Item(141, 13).Value2 = -1348  # M141: -3055 -> -1348
$ws.Cells.Item(141, 14).Value2 = -22760.0005  # N141: -22960 -> -22760.0005

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value2 = 16594.908  # H32: 16316.429 -> 16594.908
$ws.Cells.Item(32, 9).Value2 = 7482.1787  # I32: 7258.6553 -> 7482.1787
$ws.Cells.Item(32, 11).Value2 = 7482.1787  # K32: 7258.6553 -> 7482.1787
$ws.Cells.Item(32, 13).Value2 = -7195.1787  # M32: -6971.6553 -> -7195.1787

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value2 = 1394.5797  # H61: 1408.8986 -> 1394.5797
$ws.Cells.Item(61, 9).Value2 = 1394.5797  # I61: 1400.2059 -> 1394.5797
$ws.Cells.Item(61, 10).Value2 = 0  # J61: 2000 -> 0
$ws.Cells.Item(61, 11).Value2 = 1394.5797  # K61: 1400.2059 -> 1394.5797
$ws.Cells.Item(61, 12).Value2 = 0  # L61: 2000 -> 0
$ws.Cells.Item(61, 13).Value2 = -1182.5797  # M61: -1188.2059 -> -1182.5797
$ws.Cells.Item(61, 14).ClearContents() | Out-Null  # N61 was -2424

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value2 = 3403  # H74: 4192.85 -> 3403
$ws.Cells.Item(74, 9).Value2 = 1187.5454  # I74: 1292.6666 -> 1187.5454
$ws.Cells.Item(74, 10).Value2 = 6110.778  # J74: 6565.727 -> 6110.778
$ws.Cells.Item(74, 11).Value2 = 1187.5454  # K74: 1292.6666 -> 1187.5454
$ws.Cells.Item(74, 12).Value2 = 6110.778  # L74: 6565.727 -> 6110.778
$ws.Cells.Item(74, 13).Value2 = -313.5454  # M74: -418.6666 -> -313.5454
$ws.Cells.Item(74, 14).Value2 = -7858.778  # N74: -8313.726999999999 -> -7858.778

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value2 = 3403  # H77: 4192.85 -> 3403
$ws.Cells.Item(77, 9).Value2 = 1187.5454  # I77: 1292.6666 -> 1187.5454
$ws.Cells.Item(77, 10).Value2 = 6110.778  # J77: 6565.727 -> 6110.778
$ws.Cells.Item(77, 11).Value2 = 5937.727  # K77: 6463.333000000001 -> 5937.727
$ws.Cells.Item(77, 12).Value2 = 30553.89  # L77: 32828.635 -> 30553.89
$ws.Cells.Item(77, 13).Value2 = -1569.727  # M77: -2095.333000000001 -> -1569.727
$ws.Cells.Item(77, 14).Value2 = -39289.89  # N77: -41564.635 -> -39289.89

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122, 8).Value2 = 2809.923  # H122: 402248.03 -> 2809.923
$ws.Cells.Item(122, 9).Value2 = 1197.5  # I122: 557025.1 -> 1197.5
$ws.Cells.Item(122, 10).Value2 = 3526.5557  # J122: 4249.857 -> 3526.5557
$ws.Cells.Item(122, 11).Value2 = 3592.5  # K122: 1671075.3 -> 3592.5
$ws.Cells.Item(122, 12).Value2 = 10579.6671  # L122: 12749.571 -> 10579.6671
$ws.Cells.Item(122, 13).Value2 = -1142.5  # M122: -1668625.3 -> -1142.5
$ws.Cells.Item(122, 14).Value2 = -15479.6671  # N122: -17649.571 -> -15479.6671

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(125, 8).Value2 = 45000  # H125: 0 -> 45000
$ws.Cells.Item(125, 10).Value2 = 45000  # J125: 0 -> 45000
$ws.Cells.Item(125, 12).Value2 = 45000  # L125: 0 -> 45000
$ws.Cells.Item(125, 14).Value2 = -54840  # N125: None -> -54840

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value2 = 1255.9166  # H132: 1291.4572 -> 1255.9166
$ws.Cells.Item(132, 9).Value2 = 859.2059  # I132: 884.8788 -> 859.2059
$ws.Cells.Item(132, 11).Value2 = 2577.6177  # K132: 2654.6364 -> 2577.6177
$ws.Cells.Item(132, 13).Value2 = -47.61770000000024  # M132: -124.6363999999999 -> -47.61770000000024

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value2 = 1394.5797  # H136: 1408.8986 -> 1394.5797
$ws.Cells.Item(136, 9).Value2 = 1394.5797  # I136: 1400.2059 -> 1394.5797
$ws.Cells.Item(136, 10).Value2 = 0  # J136: 2000 -> 0
$ws.Cells.Item(136, 11).Value2 = 4183.7391  # K136: 4200.6177 -> 4183.7391
$ws.Cells.Item(136, 12).Value2 = 0  # L136: 6000 -> 0
$ws.Cells.Item(136, 13).Value2 = -1633.7391  # M136: -1650.6177 -> -1633.7391
$ws.Cells.Item(136, 14).ClearContents() | Out-Null  # N136 was -11100

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(139, 8).Value2 = 88500  # H139: 89000 -> 88500
$ws.Cells.Item(139, 10).Value2 = 88500  # J139: 89000 -> 88500
$ws.Cells.Item(139, 12).Value2 = 88500  # L139: 89000 -> 88500
$ws.Cells.Item(139, 14).Value2 = -98780  # N139: -99280 -> -98780

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(140, 8).Value2 = 107999.75  # H140: 108747.25 -> 107999.75
$ws.Cells.Item(140, 10).Value2 = 107999.75  # J140: 108747.25 -> 107999.75
$ws.Cells.Item(140, 12).Value2 = 107999.75  # L140: 108747.25 -> 107999.75
$ws.Cells.Item(140, 14).Value2 = -118359.75  # N140: -119107.25 -> -118359.75

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value2 = 2428.9666  # H134: 2600.9644 -> 2428.9666
$ws.Cells.Item(134, 9).Value2 = 2143.4443  # I134: 2313.24 -> 2143.4443
$ws.Cells.Item(134, 11).Value2 = 6430.3329  # K134: 6939.719999999999 -> 6430.3329
$ws.Cells.Item(134, 13).Value2 = -3895.3329  # M134: -4404.719999999999 -> -3895.3329

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value2 = 2037.52  # H31: 2389.1904 -> 2037.52
$ws.Cells.Item(31, 9).Value2 = 1029.579  # I31: 1253.1333 -> 1029.579
$ws.Cells.Item(31, 11).Value2 = 1029.579  # K31: 1253.1333 -> 1029.579
$ws.Cells.Item(31, 13).Value2 = -734.579  # M31: -958.1333 -> -734.579

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value2 = 2037.52  # H34: 2389.1904 -> 2037.52
$ws.Cells.Item(34, 9).Value2 = 1029.579  # I34: 1253.1333 -> 1029.579
$ws.Cells.Item(34, 11).Value2 = 1029.579  # K34: 1253.1333 -> 1029.579
$ws.Cells.Item(34, 13).Value2 = -827.579  # M34: -1051.1333 -> -827.579

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(58, 8).Value2 = 2956.348  # H58: 3074.818 -> 2956.348
$ws.Cells.Item(58, 9).Value2 = 1621.6364  # I58: 1748.8 -> 1621.6364
$ws.Cells.Item(58, 11).Value2 = 1621.6364  # K58: 1748.8 -> 1621.6364
$ws.Cells.Item(58, 13).Value2 = -1418.6364  # M58: -1545.8 -> -1418.6364

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(97, 8).Value2 = 0  # H97: 45000 -> 0
$ws.Cells.Item(97, 10).Value2 = 0  # J97: 45000 -> 0
$ws.Cells.Item(97, 12).Value2 = 0  # L97: 45000 -> 0
$ws.Cells.Item(97, 14).ClearContents() | Out-Null  # N97 was -46982

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(122, 8).Value2 = 6449.524  # H122: 7838.56 -> 6449.524
$ws.Cells.Item(122, 9).Value2 = 6514.091  # I122: 8293.5625 -> 6514.091
$ws.Cells.Item(122, 10).Value2 = 6378.5  # J122: 7029.6665 -> 6378.5
$ws.Cells.Item(122, 11).Value2 = 19542.273  # K122: 24880.6875 -> 19542.273
$ws.Cells.Item(122, 12).Value2 = 19135.5  # L122: 21088.9995 -> 19135.5
$ws.Cells.Item(122, 13).Value2 = -17092.273  # M122: -22430.6875 -> -17092.273
$ws.Cells.Item(122, 14).Value2 = -24035.5  # N122: -25988.9995 -> -24035.5

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(132, 8).Value2 = 4190.7856  # H132: 4311.6875 -> 4190.7856
$ws.Cells.Item(132, 10).Value2 = 6244.4  # J132: 5934 -> 6244.4
$ws.Cells.Item(132, 12).Value2 = 18733.2  # L132: 17802 -> 18733.2
$ws.Cells.Item(132, 14).Value2 = -23793.2  # N132: -22862 -> -23793.2

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value2 = 3876.0588  # H134: 3931.4375 -> 3876.0588
$ws.Cells.Item(134, 9).Value2 = 3239.3635  # I134: 3264.3 -> 3239.3635
$ws.Cells.Item(134, 11).Value2 = 9718.0905  # K134: 9792.900000000001 -> 9718.0905
$ws.Cells.Item(134, 13).Value2 = -7183.0905  # M134: -7257.900000000001 -> -7183.0905

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(136, 8).Value2 = 2956.348  # H136: 3074.818 -> 2956.348
$ws.Cells.Item(136, 9).Value2 = 1621.6364  # I136: 1748.8 -> 1621.6364
$ws.Cells.Item(136, 11).Value2 = 4864.9092  # K136: 5246.4 -> 4864.9092
$ws.Cells.Item(136, 13).Value2 = -2314.9092  # M136: -2696.4 -> -2314.9092

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(9, 8).Value2 = 1000  # H9: 0 -> 1000
$ws.Cells.Item(9, 10).Value2 = 1000  # J9: 0 -> 1000
$ws.Cells.Item(9, 12).Value2 = 3000  # L9: 0 -> 3000
$ws.Cells.Item(9, 14).Value2 = -3448  # N9: None -> -3448

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(37, 8).Value2 = 150000  # H37: 109999 -> 150000
$ws.Cells.Item(37, 10).Value2 = 150000  # J37: 109999 -> 150000
$ws.Cells.Item(37, 12).Value2 = 450000  # L37: 329997 -> 450000
$ws.Cells.Item(37, 14).Value2 = -450224  # N37: -330221 -> -450224

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(55, 8).Value2 = 102509.9  # H55: 68673 -> 102509.9
$ws.Cells.Item(55, 10).Value2 = 3949.8333  # J55: 2608.6365 -> 3949.8333
$ws.Cells.Item(55, 12).Value2 = 11849.4999  # L55: 7825.9095 -> 11849.4999
$ws.Cells.Item(55, 14).Value2 = -12203.4999  # N55: -8179.9095 -> -12203.4999

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value2 = 2208.1892  # H132: 2434.6572 -> 2208.1892
$ws.Cells.Item(132, 9).Value2 = 1703.9678  # I132: 1869.2142 -> 1703.9678
$ws.Cells.Item(132, 10).Value2 = 4813.3335  # J132: 4696.4287 -> 4813.3335
$ws.Cells.Item(132, 11).Value2 = 5111.903399999999  # K132: 5607.642599999999 -> 5111.903399999999
$ws.Cells.Item(132, 12).Value2 = 14440.0005  # L132: 14089.2861 -> 14440.0005
$ws.Cells.Item(132, 13).Value2 = -2581.903399999999  # M132: -3077.642599999999 -> -2581.903399999999
$ws.Cells.Item(132, 14).Value2 = -19500.0005  # N132: -19149.2861 -> -19500.0005

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value2 = 538.38464  # H22: 499.94116 -> 538.38464
$ws.Cells.Item(22, 9).Value2 = 568.0909  # I22: 516.6 -> 568.0909
$ws.Cells.Item(22, 11).Value2 = 568.0909  # K22: 516.6 -> 568.0909
$ws.Cells.Item(22, 13).Value2 = -273.0909  # M22: -221.6 -> -273.0909

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(27, 8).Value2 = 538.38464  # H27: 499.94116 -> 538.38464
$ws.Cells.Item(27, 9).Value2 = 568.0909  # I27: 516.6 -> 568.0909
$ws.Cells.Item(27, 11).Value2 = 568.0909  # K27: 516.6 -> 568.0909
$ws.Cells.Item(27, 13).Value2 = -461.0909  # M27: -409.6 -> -461.0909

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value2 = 1027.875  # H40: 1106.3334 -> 1027.875
$ws.Cells.Item(40, 9).Value2 = 1037.1666  # I40: 1159.5 -> 1037.1666
$ws.Cells.Item(40, 11).Value2 = 1037.1666  # K40: 1159.5 -> 1037.1666
$ws.Cells.Item(40, 13).Value2 = -901.1666  # M40: -1023.5 -> -901.1666

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(82, 8).Value2 = 1464.32  # H82: 1411.8462 -> 1464.32
$ws.Cells.Item(82, 9).Value2 = 1490.75  # I82: 1408.9412 -> 1490.75
$ws.Cells.Item(82, 11).Value2 = 1490.75  # K82: 1408.9412 -> 1490.75
$ws.Cells.Item(82, 13).Value2 = -1129.75  # M82: -1047.9412 -> -1129.75

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(85, 8).Value2 = 1464.32  # H85: 1411.8462 -> 1464.32
$ws.Cells.Item(85, 9).Value2 = 1490.75  # I85: 1408.9412 -> 1490.75
$ws.Cells.Item(85, 11).Value2 = 1490.75  # K85: 1408.9412 -> 1490.75
$ws.Cells.Item(85, 13).Value2 = -242.75  # M85: -160.9412 -> -242.75

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value2 = 4352.9165  # H132: 4523.9707 -> 4352.9165
$ws.Cells.Item(132, 9).Value2 = 3594.2083  # I132: 3789.5908 -> 3594.2083
$ws.Cells.Item(132, 11).Value2 = 10782.6249  # K132: 11368.7724 -> 10782.6249
$ws.Cells.Item(132, 13).Value2 = -8252.624899999999  # M132: -8838.7724 -> -8252.624899999999

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(60, 8).Value2 = 99997  # H60: 99999.5 -> 99997
$ws.Cells.Item(60, 9).Value2 = 99994  # I60: 0 -> 99994
$ws.Cells.Item(60, 10).Value2 = 100000  # J60: 99999.5 -> 100000
$ws.Cells.Item(60, 11).Value2 = 99994  # K60: 0 -> 99994
$ws.Cells.Item(60, 12).Value2 = 100000  # L60: 99999.5 -> 100000
$ws.Cells.Item(60, 13).Value2 = -99172  # M60: None -> -99172
$ws.Cells.Item(60, 14).Value2 = -101644  # N60: -101643.5 -> -101644

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value2 = 1356.4286  # H81: 987.0833 -> 1356.4286
$ws.Cells.Item(81, 9).Value2 = 1349.1666  # I81: 1080.625 -> 1349.1666
$ws.Cells.Item(81, 10).Value2 = 1400  # J81: 800 -> 1400
$ws.Cells.Item(81, 11).Value2 = 2698.3332  # K81: 2161.25 -> 2698.3332
$ws.Cells.Item(81, 12).Value2 = 2800  # L81: 1600 -> 2800
$ws.Cells.Item(81, 13).Value2 = -1637.3332  # M81: -1100.25 -> -1637.3332
$ws.Cells.Item(81, 14).Value2 = -4922  # N81: -3722 -> -4922

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(84, 8).Value2 = 1356.4286  # H84: 987.0833 -> 1356.4286
$ws.Cells.Item(84, 9).Value2 = 1349.1666  # I84: 1080.625 -> 1349.1666
$ws.Cells.Item(84, 10).Value2 = 1400  # J84: 800 -> 1400
$ws.Cells.Item(84, 11).Value2 = 13491.666  # K84: 10806.25 -> 13491.666
$ws.Cells.Item(84, 12).Value2 = 14000  # L84: 8000 -> 14000
$ws.Cells.Item(84, 13).Value2 = -8187.666000000001  # M84: -5502.25 -> -8187.666000000001
$ws.Cells.Item(84, 14).Value2 = -24608  # N84: -18608 -> -24608

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(126, 8).Value2 = 102372.3  # H126: 127743.25 -> 102372.3
$ws.Cells.Item(126, 9).Value2 = 143677.28  # I126: 200792.8 -> 143677.28
$ws.Cells.Item(126, 11).Value2 = 431031.84  # K126: 602378.3999999999 -> 431031.84
$ws.Cells.Item(126, 13).Value2 = -428561.84  # M126: -599908.3999999999 -> -428561.84

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value2 = 1710.2  # H132: 1565.7693 -> 1710.2
$ws.Cells.Item(132, 9).Value2 = 1637.875  # I132: 1486.909 -> 1637.875
$ws.Cells.Item(132, 11).Value2 = 4913.625  # K132: 4460.727000000001 -> 4913.625
$ws.Cells.Item(132, 13).Value2 = -2383.625  # M132: -1930.727000000001 -> -2383.625

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value2 = 68160.39999999999  # H136: 73237.07000000001 -> 68160.39999999999
$ws.Cells.Item(136, 9).Value2 = 1031.6923  # I136: 1360.4166 -> 1031.6923
$ws.Cells.Item(136, 11).Value2 = 3095.0769  # K136: 4081.2498 -> 3095.0769
$ws.Cells.Item(136, 13).Value2 = -545.0769  # M136: -1531.2498 -> -545.0769
